$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning")

# Map of short composer surnames to their full names.
$map = @{
    "Vivaldi"     = "Antonio Vivaldi"
    "Bach"        = "Jean-Sébastien Bach"
    "Beethoven"   = "Ludwig van Beethoven"
    "Chopin"      = "Frédéric Chopin"
    "Handel"      = "Georg Friedrich Haendel"
    "Mozart"      = "Wolfgang Amadeus Mozart"
    "Pachelbel"   = "Johann Pachelbel"
    "Brahms"      = "Johannes Brahms"
    "Rossini"     = "Gioachino Rossini"
    "Grieg"       = "Edvard Grieg"
    "Dvořák"      = "Antonín Dvořák"
    "Bizet"       = "Georges Bizet"
    "Schubert"    = "Franz Schubert"
    "Verdi"       = "Giuseppe Verdi"
    "Puccini"     = "Giacomo Puccini"
    "Wagner"      = "Richard Wagner"
    "Mendelssohn" = "Felix Mendelssohn Bartholdy"
    "Smetana"     = "Bedřich Smetana"
    "Liszt"       = "Franz Liszt"
    "Ravel"       = "Maurice Ravel"
    "Albinoni"    = "Tomaso Albinoni"
    "Offenbach"   = "Jacques Offenbach"
    "Borodin"     = "Alexandre Borodine"
    "Sibelius"    = "Jean Sibelius"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}

# Restore the view state: scrolled so row 15 is at top, with C15 selected.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("C15").Select()
